{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1 \u2014 Professional Summary paragraph: plain text swap, no\n// formatting change.\n//   \"...errors affecting all Black and Asian-American voters, developed...\"\n//   -> \"...errors affecting 50M voters, developed...\"\n// ---------------------------------------------------------------------\n{\n  const hits = body.search(\n    \"affecting all Black and Asian-American voters, developed geospatial ML\",\n    { matchCase: true }\n  );\n  hits.load(\"items/text\");\n  await context.sync();\n  if (hits.items.length !== 1) {\n    throw new Error(\"Summary phrase: expected 1 match, found \" + hits.items.length);\n  }\n  hits.items[0].insertText(\n    \"affecting 50M voters, developed geospatial ML\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 2 \u2014 \"Partner - Siege Analytics\" bullet: the replacement number\n// (\"50M\") becomes its own bold, colored run (matching the other\n// highlighted stats in that bullet), splitting the original single run\n// into three runs.\n// ---------------------------------------------------------------------\n{\n  const hits = body.search(\n    \"Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n    { matchCase: true }\n  );\n  hits.load(\"items/text\");\n  await context.sync();\n  if (hits.items.length !== 1) {\n    throw new Error(\"Bullet phrase: expected 1 match, found \" + hits.items.length);\n  }\n  const para = hits.items[0].paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  // Replace the whole \"affecting ... from\" span with plain new wording\n  // first (keeps the call count low and avoids touching the existing\n  // \"23%\"/\"64%\" highlighted runs that follow).\n  const span = para.search(\n    \"affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    { matchCase: true }\n  );\n  span.load(\"items/text\");\n  await context.sync();\n  if (span.items.length !== 1) {\n    throw new Error(\"Bullet span: expected 1 match, found \" + span.items.length);\n  }\n  span.items[0].insertText(\n    \"affecting 50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  // Now re-find the freshly inserted \"50M\" within the same paragraph and\n  // promote it to its own bold / colored run.\n  const num = para.search(\"50M\", { matchCase: true });\n  num.load(\"items/text\");\n  await context.sync();\n  if (num.items.length !== 1) {\n    throw new Error(\"Bullet '50M': expected 1 match, found \" + num.items.length);\n  }\n  num.items[0].font.bold = true;\n  num.items[0].font.color = \"#2C3E50\";\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 3 \u2014 \"Impact: Corrected demographic data ...\" project line:\n// plain text swap (adds \"nationwide\").\n// ---------------------------------------------------------------------\n{\n  const hits = body.search(\n    \"Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%\",\n    { matchCase: true }\n  );\n  hits.load(\"items/text\");\n  await context.sync();\n  if (hits.items.length !== 1) {\n    throw new Error(\"Impact phrase: expected 1 match, found \" + hits.items.length);\n  }\n  hits.items[0].insertText(\n    \"Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# Change 1 - Professional Summary paragraph: plain text swap, no\n# formatting change.\n#   \"...errors affecting all Black and Asian-American voters, developed...\"\n#   -> \"...errors affecting 50M voters, developed...\"\n# -----------------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$found1 = $rng1.Find.Execute(\n    \"affecting all Black and Asian-American voters, developed geospatial ML\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"affecting 50M voters, developed geospatial ML\",\n    2\n)\nif (-not $found1) {\n    throw \"Change 1: phrase not found\"\n}\n\n# -----------------------------------------------------------------------\n# Change 2 - \"Partner - Siege Analytics\" bullet: the replacement number\n# (\"50M\") becomes its own bold, colored run (matching the other\n# highlighted stats in that bullet), splitting the original single run\n# into three runs.\n# -----------------------------------------------------------------------\n$bulletPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Discovered systematic race coding errors affecting all Black and Asian-American voters*\") {\n        $bulletPara = $p\n        break\n    }\n}\nif ($null -eq $bulletPara) {\n    throw \"Change 2: bullet paragraph not found\"\n}\n\n$rng2 = $bulletPara.Range\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$found2 = $rng2.Find.Execute(\n    \"affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"affecting 50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    2\n)\nif (-not $found2) {\n    throw \"Change 2: phrase not found\"\n}\n\n# Re-find the freshly inserted \"50M\" within the same paragraph and\n# promote it to its own bold / colored run (#2C3E50, matching the other\n# highlighted stats).\n$rng3 = $bulletPara.Range\n$rng3.Find.ClearFormatting()\n$found3 = $rng3.Find.Execute(\"50M\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $found3) {\n    throw \"Change 2: '50M' not found after replace\"\n}\n$rng3.Font.Bold = 1\n$rng3.Font.Color = 5258796\n\n# -----------------------------------------------------------------------\n# Change 3 - \"Impact: Corrected demographic data ...\" project line:\n# plain text swap (adds \"nationwide\").\n# -----------------------------------------------------------------------\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Replacement.ClearFormatting()\n$found4 = $rng4.Find.Execute(\n    \"Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%\",\n    2\n)\nif (-not $found4) {\n    throw \"Change 3: phrase not found\"\n}\n"}
